$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: add an empty bordered/styled cell at M4 (continuation of header band) ---
$ws1.Range("K4").Copy()
$ws1.Range("M4").PasteSpecial(-4122)
$ws1.Range("M4").Borders.LineStyle = -4142

# Update sheet1's selection (will no longer be the active tab once Hoja2 is added/activated)
$ws1.Range("M4:M6").Select()

# --- Add the new sheet "Hoja2" right after "Hoja1" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Hoja2"
$ws2.Activate()

# --- Sheet2 data: B column steps, C/E/F columns with the difference-equation formulas ---
for ($r = 1; $r -le 8; $r++) {
    $ws2.Cells.Item($r, 2).Value = $r
}
$ws2.Range("B11").Value = 2.5
$ws2.Range("C1").Value = 6

$ws2.Range("C2").Formula = "=2*C1+2^B2-B2+3"
$ws2.Range("E2").Formula = "=2^(B2-1)*(B11+B2+4) + B2-2"
$ws2.Range("F2").Formula = "=B2+2^(B2-2)*(2*B2+13)-2"

$ws2.Range("C3").Formula = "=2*C2+2^B3-B3+3"
$ws2.Range("E3").Formula = "=2^(B3-1)*(B11+B3+4) + B3-2"
$ws2.Range("F3:F8").Formula = "=B3+2^(B3-2)*(2*B3+13)-2"

$ws2.Range("C4:C8").Formula = "=2*C3+2^B4-B4+3"
$ws2.Range("E4").Formula = "=2^(B4-1)*(B11+B4+4) + B4-2"

$ws2.Range("C10").Select()
